# Build site 2022-09-26 update for LOM3251.xlsx
#
# The underlying change removes a stray row (old row 13, which held only the
# "519033 - Carlos Yujiro Shigue" value in B/C, with no label in column A)
# which shifts every row below it up by one. The page-generation process
# that produced this workbook then also re-populated several of the shifted
# rows with different text than a pure shift would give, so those specific
# cells are corrected explicitly afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 13 (the orphan "519033 - Carlos Yujiro Shigue" row). This
#    shifts rows 14-22 up to 13-21 and updates the sheet dimension to
#    A1:C21 along with the row heights, which already land correctly.
$ws.Rows.Item(13).Delete()

# 2. Row 15 ("Programa:") needs "01/01/2012" in B/C. Copy it from B8/C8
#    (which already store that exact text as a shared string) instead of
#    typing it, so Excel doesn't reinterpret the text as a date value.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# 3. Remaining text fixups (plain text, safe to assign directly).
$ws.Range("B10").Value2 = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value2 = "519033 - Carlos Yujiro Shigue"

$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"

$ws.Range("B18").Value2 = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value2 = "519033 - Carlos Yujiro Shigue"

$ws.Range("B19").Value2 = "Provas, listas de exercícios e trabalhos práticos."
$ws.Range("C19").Value2 = "Provas, listas de exercícios e trabalhos práticos."

$ws.Range("B20").Value2 = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value2 = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

$ws.Range("B21").Value2 = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value2 = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

Write-Host "Edit applied."
